$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '88.484.85'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '3.026.38'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.32'
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '607.27'
$ws.Range("E6").Value = '  -3.79%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.358'
$ws.Range("E7").Value = '  -7.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.873'
$ws.Range("E8").Value = '  +22.52%  '
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").Value = '3.025.30'
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("E11").Value = '  +18.21%  '
$ws.Range("E12").Value = '  +3.32%  '
$ws.Range("E13").Value = '  -5.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.34'
$ws.Range("E14").Value = '  +2.10%  '
$ws.Range("D15").Value = '88.272.31'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '3.598.68'
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.42'
$ws.Range("E17").Value = '  -1.63%  '
$ws.Range("D18").Value = '3.063.64'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.38'
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000203'
$ws.Range("E20").Value = '  -2.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.21'
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '419.82'
$ws.Range("E22").Value = '  -0.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.93'
$ws.Range("E23").Value = '  +1.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.99'
$ws.Range("E24").Value = '  -3.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.36'
$ws.Range("E25").Value = '  +3.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '82.78'
$ws.Range("E26").Value = '  +5.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.50'
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").Value = '3.200.84'
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.10'
$ws.Range("E30").Value = '  +9.77%  '
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.11'
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '496.58'
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.51'
$ws.Range("E34").Value = '  -9.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.53'
$ws.Range("E35").Value = '  -2.81%  '
$ws.Range("E36").Value = '  -2.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.36'
$ws.Range("E37").Value = '  +3.67%  '
$ws.Range("E38").Value = '  -2.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.19'
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("E40").Value = '  +5.01%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  +11.22%  '
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.80'
$ws.Range("E45").Value = '  -2.70%  '
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '145.81'
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.44'
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0673'
$ws.Range("E48").Value = '  +11.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.01'
$ws.Range("E49").Value = '  +3.17%  '
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '154.54'
$ws.Range("E51").Value = '  -5.65%  '
